# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback Datetime" (column H) timestamps for the first
# handed-back file (rows 2-3) on both the zh-cn and de-de report sheets,
# as produced by regenerating the handback status report.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2:E3").Value = "2016-03-23 14:21:25"
$zhcn.Range("H2:H3").Value = "2016-03-23 14:22:03"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2:E3").Value = "2016-03-23 14:21:30"
$dede.Range("H2:H3").Value = "2016-03-23 14:22:13"
